# CUSTOM FIELDS.xlsx - "Add files via upload" edit
# Updates the three spotlight-image labels in column F to include the
# "(748x529)" recommended-dimensions hint, matching the rows whose `key`
# (column E) is spotlight_img1 / spotlight_img2 / spotlight_img3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "SECOND  SPOTLIGHT IMAGE (748x529)"
$ws.Range("F17").Value = "THIRD  SPOTLIGHT IMAGE(748x529)"
$ws.Range("F11").Value = "FIRST  SPOTLIGHT IMAGE (748x529)"

# Reflect the cursor having been left on the first edited cell.
$ws.Range("F11").Select()
